$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERTS")

# New alert rows appended to the log (rows 17-19).
$rows = @(
    @{ Row = 17; Date = "2026-02-01"; Time = "11:28:40"; Hour = "11:00"; Location = "Living Room"; Value = "CRITICAL"; Status = "FALL_DETECTED" },
    @{ Row = 18; Date = "2026-02-01"; Time = "11:28:48"; Hour = "11:00"; Location = "Living Room"; Value = "CRITICAL"; Status = "FALL_DETECTED" },
    @{ Row = 19; Date = "2026-02-01"; Time = "11:28:52"; Hour = "11:00"; Location = "Living Room"; Value = "CRITICAL"; Status = "FALL_DETECTED" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Column A holds a date-like string ("2026-02-01"). Assigning it directly
    # via .Value would be auto-parsed into a real Excel date serial (and pick
    # up a date number format), but this log keeps dates as plain text, same
    # as every other row above it. Route it through a text formula and then
    # collapse the formula down to its literal text result so the cell ends
    # up as plain text with no number-format override.
    $dateCell = $ws.Cells.Item($rowIndex, 1)
    $dateCell.Formula = '="' + $r.Date + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    $ws.Cells.Item($rowIndex, 2).Value = $r.Time
    $ws.Cells.Item($rowIndex, 3).Value = $r.Hour
    $ws.Cells.Item($rowIndex, 4).Value = $r.Location
    $ws.Cells.Item($rowIndex, 5).Value = $r.Value
    $ws.Cells.Item($rowIndex, 6).Value = $r.Status
}

$excel.CutCopyMode = 0
